# Generate Report for Handback
# Refresh the localization status report: the handback status moves from
# "Ready for handoff" to "Handed back: in sync with en-US", the per-language
# "Latest Handback DateTime" stamps advance, and the now-resolved
# "Error Detail" (stale handback-version warning) is cleared out.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns (E2, F2)
$overview.Range("E2").Value = $status
$overview.Range("F2").Value = $status

# zh-cn detail sheet
$zhcn.Range("C2").Value = $status
$zhcn.Range("K2").Value = "2016-08-22 20:48:52"
$zhcn.Range("P2").Value = ""

# de-de detail sheet
$dede.Range("C2").Value = $status
$dede.Range("K2").Value = "2016-08-22 20:49:00"
$dede.Range("P2").Value = ""

# Re-fit the columns whose text width changed after the refresh above.
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null
$zhcn.Columns.Item(3).AutoFit() | Out-Null
$zhcn.Columns.Item(16).AutoFit() | Out-Null
$dede.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(16).AutoFit() | Out-Null
